# --- Countries & provincias Spain data refresh (25 May 2020, 14:35 -> 15:05) ---
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Last updated" timestamp banner in A1
$ws.Range("A1").Value2 = "Datos actualizados a 25 de Mayo de 2020 a las 15:05"

# Row 4: Estados Unidos
$ws.Range("B4").Value2 = 1688290
$ws.Range("C4").Value2 = 1854
$ws.Range("E4").Value2 = 1137197
$ws.Range("G4").Value2 = 48
$ws.Range("H4").Value2 = 99348

# Row 11: Alemania
$ws.Range("B11").Value2 = 180331
$ws.Range("C11").Value2 = 3
$ws.Range("E11").Value2 = 10760

# Row 13: India
$ws.Range("B13").Value2 = 140146
$ws.Range("C13").Value2 = 1610
$ws.Range("D13").Value2 = 58174
$ws.Range("E13").Value2 = 77931
$ws.Range("G13").Value2 = 17
$ws.Range("H13").Value2 = 4041

# Row 18: Arabia Saudita
$ws.Range("B18").Value2 = 74795
$ws.Range("C18").Value2 = 2235
$ws.Range("D18").Value2 = 45668
$ws.Range("E18").Value2 = 28728
$ws.Range("G18").Value2 = 9
$ws.Range("H18").Value2 = 399

# Row 47: Argentina
$ws.Range("D47").Value2 = 3999
$ws.Range("E47").Value2 = 7621
$ws.Range("G47").Value2 = 4
$ws.Range("H47").Value2 = 456

# Row 50: Serbia
$ws.Range("A50").Value2 = "Serbia"
$ws.Range("B50").Value2 = 11193
$ws.Range("C50").Value2 = 34
$ws.Range("D50").Value2 = 5920
$ws.Range("E50").Value2 = 5034
$ws.Range("H50").Value2 = 239

# Row 51: Afganistan
$ws.Range("A51").Value2 = "Afganistan"
$ws.Range("B51").Value2 = 11173
$ws.Range("C51").Value2 = 591
$ws.Range("D51").Value2 = 1097
$ws.Range("E51").Value2 = 9857
$ws.Range("G51").Value2 = 1
$ws.Range("H51").Value2 = 219

# Row 53: Barein
$ws.Range("B53").Value2 = 9164
$ws.Range("C53").Value2 = 26
$ws.Range("D53").Value2 = 4753
$ws.Range("E53").Value2 = 4397

# Row 108: Mali
$ws.Range("A108").Value2 = "Mali"
$ws.Range("B108").Value2 = 1059
$ws.Range("C108").Value2 = 29
$ws.Range("D108").Value2 = 604
$ws.Range("E108").Value2 = 388
$ws.Range("G108").Value2 = 2
$ws.Range("H108").Value2 = 67

# Row 109: Tunez
$ws.Range("A109").Value2 = "Tunez"
$ws.Range("B109").Value2 = 1051
$ws.Range("C109").Value2 = 3
$ws.Range("D109").Value2 = 917
$ws.Range("E109").Value2 = 86
$ws.Range("H109").Value2 = 48

# Row 110: Letonia
$ws.Range("A110").Value2 = "Letonia"
$ws.Range("B110").Value2 = 1049
$ws.Range("C110").Value2 = 2
$ws.Range("D110").Value2 = 712
$ws.Range("E110").Value2 = 315
$ws.Range("G110").Value2 = 0
$ws.Range("H110").Value2 = 22

# Row 111: Guinea Ecuatorial
$ws.Range("A111").Value2 = "Guinea Ecuatorial"
$ws.Range("B111").Value2 = 1043
$ws.Range("C111").Value2 = 83
$ws.Range("D111").Value2 = 165
$ws.Range("E111").Value2 = 866
$ws.Range("G111").Value2 = 1
$ws.Range("H111").Value2 = 12

# Row 122: Sierra Leona
$ws.Range("A122").Value2 = "Sierra Leona"
$ws.Range("B122").Value2 = 735
$ws.Range("C122").Value2 = 28
$ws.Range("D122").Value2 = 293
$ws.Range("E122").Value2 = 400
$ws.Range("G122").Value2 = 2
$ws.Range("H122").Value2 = 42

# Row 123: Georgia
$ws.Range("A123").Value2 = "Georgia"
$ws.Range("B123").Value2 = 731
$ws.Range("C123").Value2 = 1
$ws.Range("D123").Value2 = 526
$ws.Range("E123").Value2 = 193
$ws.Range("H123").Value2 = 12

# Row 124: Crucero
$ws.Range("A124").Value2 = "Crucero"
$ws.Range("B124").Value2 = 712
$ws.Range("D124").Value2 = 651
$ws.Range("E124").Value2 = 48
$ws.Range("H124").Value2 = 13

# Row 125: Jordania
$ws.Range("A125").Value2 = "Jordania"
$ws.Range("B125").Value2 = 708
$ws.Range("D125").Value2 = 471
$ws.Range("E125").Value2 = 228
$ws.Range("H125").Value2 = 9

# Row 140: Cabo Verde
$ws.Range("A140").Value2 = "Cabo Verde"
$ws.Range("B140").Value2 = 390
$ws.Range("C140").Value2 = 10
$ws.Range("D140").Value2 = 155
$ws.Range("E140").Value2 = 232
$ws.Range("H140").Value2 = 3

# Row 141: Togo
$ws.Range("A141").Value2 = "Togo"
$ws.Range("B141").Value2 = 381
$ws.Range("D141").Value2 = 141
$ws.Range("E141").Value2 = 228
$ws.Range("H141").Value2 = 12
